# Kentucky overview workbook: convert numeric "count" cells to text cells
# (same digits, just stored as text) and backfill the six previously-zero
# counties in the County sheet with text placeholder values, adding the
# new statewide "Total" row at the bottom of that sheet.

function Set-TextValue {
    param($cell, [string]$text)
    # Force text storage: without this, values that look numeric (plain
    # integers, "$0", "0.00%", "1,138", ...) get auto-coerced back into a
    # number by Excel's input parser.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    # Drop back to the Normal style so we don't leave a stray cell style
    # behind just because we touched NumberFormat.
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overall": A2 1138 -> "1,138"
# ---------------------------------------------------------------------
$wsOverall = $wb.Worksheets.Item("Overall")
Set-TextValue $wsOverall.Cells.Item(2, 1) "1,138"

# ---------------------------------------------------------------------
# Sheet "County": B2:B104 numeric -> text (same digits); the six
# zero-filer counties (rows 105-110) get placeholder text values; add a
# new "Total" row 111.
# ---------------------------------------------------------------------
$wsCounty = $wb.Worksheets.Item("County")

$countyCounts = @{2="1";3="4";4="1";5="1";6="8";7="7";8="11";9="4";10="15";11="14";12="1";13="4";14="2";15="8";16="3";17="2";18="6";19="21";20="3";21="2";22="15";23="8";24="3";25="2";26="1";27="3";28="37";29="2";30="2";31="150";32="1";33="11";34="29";35="1";36="1";37="1";38="2";39="10";40="1";41="2";42="21";43="5";44="1";45="2";46="18";47="2";48="14";49="1";50="292";51="9";52="5";53="45";54="4";55="3";56="1";57="16";58="1";59="2";60="4";61="4";62="2";63="3";64="3";65="1";66="20";67="1";68="5";69="5";70="8";71="23";72="1";73="3";74="5";75="1";76="1";77="10";78="2";79="2";80="10";81="2";82="9";83="5";84="2";85="19";86="14";87="1";88="12";89="2";90="8";91="6";92="12";93="8";94="9";95="3";96="2";97="3";98="5";99="44";100="1";101="2";102="7";103="1";104="5"}

foreach ($row in $countyCounts.Keys) {
    Set-TextValue $wsCounty.Cells.Item([int]$row, 2) $countyCounts[$row]
}

# Rows 105-110: the six counties with no 990 filers w/ gov grants.
$zeroCountyRows = 105,106,107,108,109,110
foreach ($row in $zeroCountyRows) {
    Set-TextValue $wsCounty.Cells.Item($row, 2) "0.00%"
    Set-TextValue $wsCounty.Cells.Item($row, 3) "`$0"
    Set-TextValue $wsCounty.Cells.Item($row, 4) "0.00%"
    Set-TextValue $wsCounty.Cells.Item($row, 5) "0.00%"
    Set-TextValue $wsCounty.Cells.Item($row, 6) "0.00%"
}

# New statewide Total row at the bottom of the County sheet.
Set-TextValue $wsCounty.Cells.Item(111, 1) "Total"
Set-TextValue $wsCounty.Cells.Item(111, 2) "1,138"
Set-TextValue $wsCounty.Cells.Item(111, 3) "`$2,425,907,719"
Set-TextValue $wsCounty.Cells.Item(111, 4) "9.25%"
Set-TextValue $wsCounty.Cells.Item(111, 5) "-10.73%"
Set-TextValue $wsCounty.Cells.Item(111, 6) "65.64%"

# ---------------------------------------------------------------------
# Sheet "Congressional District": B2:B7 numeric -> text; Total (B8) ->
# "1,138"
# ---------------------------------------------------------------------
$wsCongress = $wb.Worksheets.Item("Congressional District")

$congressCounts = @{2="196";3="150";4="288";5="118";6="157";7="229"}
foreach ($row in $congressCounts.Keys) {
    Set-TextValue $wsCongress.Cells.Item([int]$row, 2) $congressCounts[$row]
}
Set-TextValue $wsCongress.Cells.Item(8, 2) "1,138"

# ---------------------------------------------------------------------
# Sheet "Size": B2:B7 numeric -> text; Total (B8) -> "1,138"
# ---------------------------------------------------------------------
$wsSize = $wb.Worksheets.Item("Size")

$sizeCounts = @{2="376";3="268";4="175";5="76";6="172";7="71"}
foreach ($row in $sizeCounts.Keys) {
    Set-TextValue $wsSize.Cells.Item([int]$row, 2) $sizeCounts[$row]
}
Set-TextValue $wsSize.Cells.Item(8, 2) "1,138"

# ---------------------------------------------------------------------
# Sheet "Subsector": B2:B12 numeric -> text; Total (B13) -> "1,138"
# ---------------------------------------------------------------------
$wsSubsector = $wb.Worksheets.Item("Subsector")

$subsectorCounts = @{2="100";3="91";4="43";5="123";6="33";7="358";8="10";9="84";10="26";11="255";12="15"}
foreach ($row in $subsectorCounts.Keys) {
    Set-TextValue $wsSubsector.Cells.Item([int]$row, 2) $subsectorCounts[$row]
}
Set-TextValue $wsSubsector.Cells.Item(13, 2) "1,138"
